$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.095.44"
$ws.Range("E2").Value = "  +1.77%  "

$ws.Range("D3").Value = "2.254.13"
$ws.Range("E3").Value = "  +1.00%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "272.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0931"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.31%  "

$ws.Range("D14").Value = "2.590.63"
$ws.Range("E14").Value = "  +0.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.64%  "

$ws.Range("D16").Value = "2.271.97"
$ws.Range("E16").Value = "  +1.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.801"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.01%  "

$ws.Range("D18").Value = "44.018.50"
$ws.Range("E18").Value = "  +1.70%  "

$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("E20").Value = "  -0.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.11%  "

$ws.Range("E24").Value = "  -4.48%  "

$ws.Range("E26").Value = "  +13.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.35%  "

$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.23%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0897"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.73%  "

$ws.Range("E35").Value = "  +1.34%  "

$ws.Range("E36").Value = "  +3.62%  "

$ws.Range("E37").Value = "  -4.35%  "

$ws.Range("E38").Value = "  -1.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +21.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.12%  "

$ws.Range("E41").Value = "  +2.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.04%  "

$ws.Range("E44").Value = "  +0.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.93%  "

$ws.Range("E46").Value = "  +0.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.54%  "

$ws.Range("E48").Value = "  +4.84%  "

$ws.Range("E49").Value = "  +1.71%  "

$ws.Range("E50").Value = "  +1.52%  "

$ws.Range("E51").Value = "  -8.39%  "
